$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 10000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10936
$ws.Range("M21").ClearContents()
# Row 23
$ws.Range("H23").Value = 10000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 10000
$ws.Range("N23").Value = -10468
$ws.Range("M23").ClearContents()
# Row 29
$ws.Range("H29").Value = 486
$ws.Range("J29").Value = 472
$ws.Range("L29").Value = 1416
$ws.Range("N29").Value = -1978
# Row 38
$ws.Range("H38").Value = 137.6
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
# Row 58
$ws.Range("H58").Value = 755.6667
$ws.Range("I58").Value = 125
$ws.Range("J58").Value = 2017
$ws.Range("K58").Value = 375
$ws.Range("L58").Value = 6051
$ws.Range("M58").Value = -225
$ws.Range("N58").Value = -6351
# Row 87
$ws.Range("H87").Value = 38638
$ws.Range("J87").Value = 38638
$ws.Range("L87").Value = 38638
$ws.Range("N87").Value = -41134
# Row 90
$ws.Range("H90").Value = 38638
$ws.Range("J90").Value = 38638
$ws.Range("L90").Value = 115914
$ws.Range("N90").Value = -128394
# Row 97
$ws.Range("H97").Value = 2925
$ws.Range("I97").Value = 540
$ws.Range("J97").Value = 5310
$ws.Range("K97").Value = 1620
$ws.Range("L97").Value = 15930
$ws.Range("M97").Value = -1124
$ws.Range("N97").Value = -16922
# Row 100
$ws.Range("H100").Value = 3000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 3000
$ws.Range("N100").Value = -4082
$ws.Range("M100").ClearContents()
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 132
$ws.Range("H132").Value = 1090735.8
$ws.Range("I132").Value = 1943.0714
$ws.Range("J132").Value = 16333833
$ws.Range("K132").Value = 5829.2142
$ws.Range("L132").Value = 49001499
$ws.Range("M132").Value = -3299.2142
$ws.Range("N132").Value = -49006559
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
# Row 137
$ws.Range("H137").Value = 6001745.5
$ws.Range("I137").Value = 15001094
$ws.Range("K137").Value = 45003282
$ws.Range("M137").Value = -45000732
# Row 139
$ws.Range("H139").Value = 41868
$ws.Range("J139").Value = 41868
$ws.Range("L139").Value = 41868
$ws.Range("N139").Value = -52148

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 66801628
$ws.Range("I61").Value = 100101570
$ws.Range("J61").Value = 201742.8
$ws.Range("K61").Value = 100101570
$ws.Range("L61").Value = 201742.8
$ws.Range("M61").Value = -100101358
$ws.Range("N61").Value = -202166.8
# Row 102
$ws.Range("H102").Value = 76428570
$ws.Range("I102").Value = 76428570
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 76428570
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -76426948
$ws.Range("N102").ClearContents()
# Row 136
$ws.Range("H136").Value = 66801628
$ws.Range("I136").Value = 100101570
$ws.Range("J136").Value = 201742.8
$ws.Range("K136").Value = 300304710
$ws.Range("L136").Value = 605228.3999999999
$ws.Range("M136").Value = -300302160
$ws.Range("N136").Value = -610328.3999999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 9420.467000000001
$ws.Range("I86").Value = 18524
$ws.Range("J86").Value = 2458.9412
$ws.Range("K86").Value = 18524
$ws.Range("L86").Value = 2458.9412
$ws.Range("M86").Value = -17401
$ws.Range("N86").Value = -4704.9412
# Row 89
$ws.Range("H89").Value = 9420.467000000001
$ws.Range("I89").Value = 18524
$ws.Range("J89").Value = 2458.9412
$ws.Range("K89").Value = 92620
$ws.Range("L89").Value = 12294.706
$ws.Range("M89").Value = -87004
$ws.Range("N89").Value = -23526.706
# Row 134
$ws.Range("H134").Value = 7234
$ws.Range("I134").Value = 6549.8823
$ws.Range("K134").Value = 19649.6469
$ws.Range("M134").Value = -17114.6469

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 901.02856
$ws.Range("J131").Value = 980.5517
$ws.Range("L131").Value = 2941.6551
$ws.Range("N131").Value = -13021.6551

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 57913.633
$ws.Range("I70").Value = 81907.69500000001
$ws.Range("J70").Value = 5926.5
$ws.Range("K70").Value = 81907.69500000001
$ws.Range("L70").Value = 5926.5
$ws.Range("M70").Value = -81637.69500000001
$ws.Range("N70").Value = -6466.5
# Row 73
$ws.Range("H73").Value = 57913.633
$ws.Range("I73").Value = 81907.69500000001
$ws.Range("J73").Value = 5926.5
$ws.Range("K73").Value = 81907.69500000001
$ws.Range("L73").Value = 5926.5
$ws.Range("M73").Value = -80971.69500000001
$ws.Range("N73").Value = -7798.5
# Row 122
$ws.Range("H122").Value = 2082.182
$ws.Range("I122").Value = 1697.5
$ws.Range("J122").Value = 2302
$ws.Range("K122").Value = 5092.5
$ws.Range("L122").Value = 6906
$ws.Range("M122").Value = -2642.5
$ws.Range("N122").Value = -11806
# Row 126
$ws.Range("H126").Value = 3304.8
$ws.Range("I126").Value = 1674.6666
$ws.Range("J126").Value = 5750
$ws.Range("K126").Value = 5023.9998
$ws.Range("L126").Value = 17250
$ws.Range("M126").Value = -2553.9998
$ws.Range("N126").Value = -22190
# Row 132
$ws.Range("H132").Value = 48924.812
$ws.Range("I132").Value = 47876.684
$ws.Range("J132").Value = 50022.855
$ws.Range("K132").Value = 143630.052
$ws.Range("L132").Value = 150068.565
$ws.Range("M132").Value = -141100.052
$ws.Range("N132").Value = -155128.565

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 42908.16
$ws.Range("I132").Value = 2527.111
$ws.Range("J132").Value = 65622.5
$ws.Range("K132").Value = 7581.333
$ws.Range("L132").Value = 196867.5
$ws.Range("M132").Value = -5051.333
$ws.Range("N132").Value = -201927.5
# Row 136
$ws.Range("H136").Value = 74456.27
$ws.Range("I136").Value = 38616.332
$ws.Range("J136").Value = 148893.08
$ws.Range("K136").Value = 115848.996
$ws.Range("L136").Value = 446679.24
$ws.Range("M136").Value = -113298.996
$ws.Range("N136").Value = -451779.24

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 76188.7
$ws.Range("I132").Value = 60664.59
$ws.Range("J132").Value = 102579.7
$ws.Range("K132").Value = 181993.77
$ws.Range("L132").Value = 307739.1
$ws.Range("M132").Value = -179463.77
$ws.Range("N132").Value = -312799.1
# Row 136
$ws.Range("H136").Value = 54861.164
$ws.Range("I136").Value = 35276.38
$ws.Range("J136").Value = 125856
$ws.Range("K136").Value = 105829.14
$ws.Range("L136").Value = 377568
$ws.Range("M136").Value = -103279.14
$ws.Range("N136").Value = -382668
